# Insert a new price-record row at row 453 (Región de La Araucanía,
# $/docena de atados (2 kilos), 50 units @ 6000) and push the existing
# rows 453:468 down to 454:469, matching the weekly Fruta/Hortaliza
# refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a fresh row above the old row 453.
$ws.Range("A453").EntireRow.Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A453").Value = 4
$ws.Range("B453").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C453").Value = "Los Lagos"
$ws.Range("D453").Value = 45075
$ws.Range("E453").Value = 10
$ws.Range("F453").Value = 100112040
$ws.Range("G453").Value = "Cilantro"
$ws.Range("H453").Value = "Sin especificar"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 50
$ws.Range("K453").Value = 6000
$ws.Range("L453").Value = 6000
$ws.Range("M453").Value = 6000
$ws.Range("N453").Value = "$/docena de atados (2 kilos)"
$ws.Range("O453").Value = "Región de La Araucanía"
$ws.Range("P453").Value = 3000
$ws.Range("Q453").Value = 2
$ws.Range("R453").Value = "Hortaliza"
